$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 5555
$ws.Range("I40").Value = 5555
$ws.Range("K40").Value = 5555
$ws.Range("M40").Value = -5380
$ws.Range("H64").Value = 1973.3334
$ws.Range("I64").Value = 1973.3334
$ws.Range("K64").Value = 1973.3334
$ws.Range("M64").Value = -1725.3334
$ws.Range("H67").Value = 1973.3334
$ws.Range("I67").Value = 1973.3334
$ws.Range("K67").Value = 1973.3334
$ws.Range("M67").Value = -1115.3334
$ws.Range("H70").Value = 2524.923
$ws.Range("J70").Value = 3404.2856
$ws.Range("L70").Value = 10212.8568
$ws.Range("N70").Value = -10752.8568
$ws.Range("H73").Value = 2524.923
$ws.Range("J73").Value = 3404.2856
$ws.Range("L73").Value = 10212.8568
$ws.Range("N73").Value = -12084.8568
$ws.Range("H74").Value = 500000
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("H76").Value = 5792.875
$ws.Range("I76").Value = 5070.8
$ws.Range("J76").Value = 6996.3335
$ws.Range("K76").Value = 5070.8
$ws.Range("L76").Value = 6996.3335
$ws.Range("M76").Value = -4755.8
$ws.Range("N76").Value = -7626.3335
$ws.Range("H77").Value = 500000
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("H79").Value = 5792.875
$ws.Range("I79").Value = 5070.8
$ws.Range("J79").Value = 6996.3335
$ws.Range("K79").Value = 5070.8
$ws.Range("L79").Value = 6996.3335
$ws.Range("M79").Value = -3978.8
$ws.Range("N79").Value = -9180.333500000001
$ws.Range("H100").Value = 875
$ws.Range("I100").Value = 875
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 875
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -334
$ws.Range("N100").ClearContents()
$ws.Range("H111").Value = 1092.5714
$ws.Range("I111").Value = 1092.5714
$ws.Range("K111").Value = 3277.7142
$ws.Range("M111").Value = -210.7142000000003
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H131").Value = 1356.2858
$ws.Range("I131").Value = 1173.5
$ws.Range("J131").Value = 1600
$ws.Range("K131").Value = 3520.5
$ws.Range("L131").Value = 4800
$ws.Range("M131").Value = 1519.5
$ws.Range("N131").Value = -14880
$ws.Range("H132").Value = 5403.857
$ws.Range("I132").Value = 5306.1665
$ws.Range("K132").Value = 15918.4995
$ws.Range("M132").Value = -13388.4995
$ws.Range("H141").Value = 4996
$ws.Range("I141").Value = 4701.8237
$ws.Range("K141").Value = 14105.4711
$ws.Range("M141").Value = -8925.471099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 215.57143
$ws.Range("J2").Value = 247.5
$ws.Range("L2").Value = 247.5
$ws.Range("N2").Value = -473.5
$ws.Range("H32").Value = 8391.936
$ws.Range("I32").Value = 7005
$ws.Range("K32").Value = 7005
$ws.Range("M32").Value = -6718
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()
$ws.Range("H63").Value = 2326
$ws.Range("I63").Value = 1361.6
$ws.Range("K63").Value = 1361.6
$ws.Range("M63").Value = -675.5999999999999
$ws.Range("H66").Value = 2326
$ws.Range("I66").Value = 1361.6
$ws.Range("K66").Value = 6808
$ws.Range("M66").Value = -3376
$ws.Range("H97").Value = 702.44446
$ws.Range("I97").Value = 665.875
$ws.Range("K97").Value = 665.875
$ws.Range("M97").Value = -169.875
$ws.Range("H116").Value = 215.57143
$ws.Range("J116").Value = 247.5
$ws.Range("L116").Value = 247.5
$ws.Range("N116").Value = -4835.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 215.57143
$ws.Range("J3").Value = 247.5
$ws.Range("L3").Value = 247.5
$ws.Range("N3").Value = -475.5
$ws.Range("H54").Value = 37500
$ws.Range("I54").Value = 37500
$ws.Range("K54").Value = 37500
$ws.Range("M54").Value = -37016
$ws.Range("H86").Value = 3469.077
$ws.Range("I86").Value = 3199.889
$ws.Range("K86").Value = 3199.889
$ws.Range("M86").Value = -2076.889
$ws.Range("H89").Value = 3469.077
$ws.Range("I89").Value = 3199.889
$ws.Range("K89").Value = 15999.445
$ws.Range("M89").Value = -10383.445

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4331.647
$ws.Range("J16").Value = 7560.7144
$ws.Range("L16").Value = 7560.7144
$ws.Range("N16").Value = -8134.7144
$ws.Range("H50").Value = 30000
$ws.Range("J50").Value = 30000
$ws.Range("L50").Value = 30000
$ws.Range("N50").Value = -31250
$ws.Range("H60").Value = 20160.666
$ws.Range("J60").Value = 24994.285
$ws.Range("L60").Value = 24994.285
$ws.Range("N60").Value = -26016.285
$ws.Range("H62").Value = 6622.5
$ws.Range("I62").Value = 6622.5
$ws.Range("K62").Value = 6622.5
$ws.Range("M62").Value = -5998.5
$ws.Range("H65").Value = 6622.5
$ws.Range("I65").Value = 6622.5
$ws.Range("K65").Value = 33112.5
$ws.Range("M65").Value = -29992.5
$ws.Range("H113").Value = 4331.647
$ws.Range("J113").Value = 7560.7144
$ws.Range("L113").Value = 7560.7144
$ws.Range("N113").Value = -11900.7144
$ws.Range("H132").Value = 1997.7
$ws.Range("I132").Value = 1997.7
$ws.Range("K132").Value = 5993.1
$ws.Range("M132").Value = -3463.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 97936
$ws.Range("J37").Value = 97936
$ws.Range("L37").Value = 293808
$ws.Range("N37").Value = -294032

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4271.25
$ws.Range("I80").Value = 1835
$ws.Range("K80").Value = 1835
$ws.Range("M80").Value = -837
$ws.Range("H83").Value = 4271.25
$ws.Range("I83").Value = 1835
$ws.Range("K83").Value = 9175
$ws.Range("M83").Value = -4183
$ws.Range("H97").Value = 776.3158
$ws.Range("I97").Value = 791.17645
$ws.Range("J97").Value = 650
$ws.Range("K97").Value = 791.17645
$ws.Range("L97").Value = 650
$ws.Range("M97").Value = -295.17645
$ws.Range("N97").Value = -1642
$ws.Range("H122").Value = 125000000
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 1392.3334
$ws.Range("I132").Value = 1253.875
$ws.Range("K132").Value = 3761.625
$ws.Range("M132").Value = -1231.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3197.4
$ws.Range("I46").Value = 3995.5
$ws.Range("J46").Value = 2665.3333
$ws.Range("K46").Value = 3995.5
$ws.Range("L46").Value = 2665.3333
$ws.Range("M46").Value = -3807.5
$ws.Range("N46").Value = -3041.3333
$ws.Range("H82").Value = 1809.091
$ws.Range("J82").Value = 1816.6666
$ws.Range("L82").Value = 1816.6666
$ws.Range("N82").Value = -2538.6666
$ws.Range("H85").Value = 1809.091
$ws.Range("J85").Value = 1816.6666
$ws.Range("L85").Value = 1816.6666
$ws.Range("N85").Value = -4312.6666
$ws.Range("H93").Value = 221.25
$ws.Range("I93").Value = 161.66667
$ws.Range("K93").Value = 161.66667
$ws.Range("M93").Value = 1086.33333
$ws.Range("H106").Value = 11497.5
$ws.Range("J106").Value = 11497.5
$ws.Range("L106").Value = 11497.5
$ws.Range("N106").Value = -14021.5
$ws.Range("H122").Value = 3977.923
$ws.Range("I122").Value = 3610.2727
$ws.Range("K122").Value = 10830.8181
$ws.Range("M122").Value = -8380.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 874.2
$ws.Range("I122").Value = 874.2
$ws.Range("K122").Value = 2622.6
$ws.Range("M122").Value = -172.6000000000004
$ws.Range("H126").Value = 4419.2
$ws.Range("I126").Value = 4419.2
$ws.Range("K126").Value = 13257.6
$ws.Range("M126").Value = -10787.6
